# Auto-applies cell value updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainCell($cell, $value) {
    $ws.Range($cell).Value = $value
}

Set-TextCell "D2" "330.91"
Set-TextCell "E2" "0.53%"
Set-TextCell "E3" "0.34%"
Set-TextCell "D4" "5.699"
Set-TextCell "E4" "-0.17%"
Set-TextCell "D5" "0.08429"
Set-TextCell "E5" "4.55%"
Set-TextCell "E6" "1.11%"
Set-TextCell "E7" "-0.42%"
Set-TextCell "D8" "1.982"
Set-TextCell "E8" "-3.48%"
Set-TextCell "D10" "0.9270"
Set-TextCell "E10" "0.66%"
Set-TextCell "D11" "0.1254"
Set-TextCell "E11" "0.63%"
Set-TextCell "D12" "0.1978"
Set-TextCell "E12" "1.68%"
Set-TextCell "D13" "0.09571"
Set-TextCell "E13" "3.03%"
Set-TextCell "D14" "0.03965"
Set-TextCell "E14" "8.17%"
Set-TextCell "E15" "0.87%"
Set-TextCell "D16" "0.001304"
Set-TextCell "E16" "0.41%"
Set-PlainCell "B17" "TigerCash"
Set-PlainCell "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D17" "0.006115"
Set-TextCell "E17" "-1.07%"
Set-PlainCell "B18" "LEO"
Set-PlainCell "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D18" "3.436"
Set-TextCell "E18" "1.60%"
Set-PlainCell "B19" "BitpandaEcosystemToken"
Set-PlainCell "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D19" "0.3511"
Set-TextCell "E19" "0.83%"
Set-PlainCell "B20" "MCDex"
Set-PlainCell "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D20" "9.167"
Set-TextCell "E20" "10.57%"
Set-PlainCell "B21" "ProBitToken"
Set-PlainCell "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell "D21" "0.1364"
Set-TextCell "E21" "-3.63%"
Set-PlainCell "B22" "ZBToken"
Set-PlainCell "C22" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell "D22" "0.2512"
Set-TextCell "E22" "-5.24%"
Set-PlainCell "B23" "CoinExToken"
Set-PlainCell "C23" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D23" "0.04403"
Set-TextCell "E23" "-0.59%"
Set-TextCell "D24" "0.001247"
Set-TextCell "E24" "-1.08%"
Set-TextCell "D25" "0.004371"
Set-TextCell "E25" "0.75%"
Set-TextCell "E26" "-3.95%"
Set-TextCell "E27" "0.08%"
Set-TextCell "D39" "0.02831"
Set-TextCell "E39" "0.30%"
Set-TextCell "D40" "0.05518"
Set-TextCell "E40" "0.93%"
Set-TextCell "D41" "0.007919"
Set-TextCell "E41" "3.94%"
Set-TextCell "D42" "0.1439"
Set-TextCell "E42" "1.61%"
Set-TextCell "D43" "0.008957"
Set-TextCell "E43" "-9.90%"
Set-TextCell "D44" "0.002094"
Set-TextCell "E44" "-0.88%"
Set-TextCell "E45" "-7.62%"
Set-TextCell "D46" "0.00007328"
Set-TextCell "E46" "8.98%"
Set-TextCell "E47" "0.22%"
Set-TextCell "D48" "0.003212"
Set-TextCell "E48" "1.24%"
Set-TextCell "D49" "0.002282"
Set-TextCell "E49" "0.15%"
Set-TextCell "D50" "0.00002104"
Set-TextCell "E50" "0.22%"
Set-TextCell "E51" "0.22%"
